$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two oldest entries at the bottom of the table
# (미래에셋비전스팩5호 row 20, 한국스팩14호 row 21)
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(20).Delete()

# Insert two new SPAC book-building entries right after row 2 (아이빔테크놀로지)
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = "NH스팩31호"
$ws.Range("B3").Value = "2024.07.09~07.10"
$ws.Range("C3").Value = "2,000~2,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 12000
$ws.Range("F3").Value = "NH투자증권"

$ws.Range("A4").Value = "SK증권스팩13호"
$ws.Range("B4").Value = "2024.07.09~07.10"
$ws.Range("C4").Value = "2,000~2,000"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = 8000
$ws.Range("F4").Value = "SK증권"

# 하이젠알앤엠 (now row 17 after the two inserts above) has its final
# offering price set: 확정공모가 goes from "-" to 7000
$ws.Range("D17").Value = 7000
